# Add a "Deserialization" bullet (level-2 indent) right after "Extracting
# values from responses" in the capstone-assignment slide's content
# placeholder.

$p = $ppt.ActivePresentation

# Locate the slide/shape that contains the "Extracting values from
# responses" bullet, so the script is resilient to slide re-numbering.
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -like "*Extracting values from responses*") {
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) {
        break
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph whose text is exactly "Extracting values from responses".
$paraCount = $tr.Paragraphs().Count
$anchorIndex = -1
for ($k = 1; $k -le $paraCount; $k++) {
    $para = $tr.Paragraphs($k, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq "Extracting values from responses") {
        $anchorIndex = $k
        break
    }
}

$anchorPara = $tr.Paragraphs($anchorIndex, 1)

# Insert a new paragraph right after it. The leading carriage return creates
# the new paragraph break so "Deserialization" lands in its own <a:p>,
# inheriting the level-1 bullet formatting (green Courier New, "_" bullet
# char) of the anchor paragraph automatically.
$inserted = $anchorPara.InsertAfter("`rDeserialization")
